$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right
#    after the title heading (paragraph 2).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Arthur's Fortune for Free: Review
#    & Bonuses") right before the last paragraph in the document (the one
#    that currently holds the feature-image prompt text), mirroring the
#    structure of the old meta-description paragraph's bold run.
$n = $d.Paragraphs.Count
$priorPara = $d.Paragraphs.Item($n - 1)
$insertPoint = $priorPara.Range.Duplicate
$insertPoint.Collapse(0)

$titleOoxml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Arthur''s Fortune for Free: Review &amp; Bonuses</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($titleOoxml)

# 3. Replace the text of the final (italic) paragraph with the review
#    blurb that used to live in the meta-description paragraph, keeping
#    its existing italic character formatting intact.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bodyRange = $lastPara.Range.Duplicate
$bodyRange.MoveEnd(1, -1) | Out-Null
$bodyRange.Text = "Read our review of Arthur's Fortune, a medieval themed slot with free spins, random conversion feature, and bonus feature. Play for free and win big."
